# Updated symbol list on Fri Dec 23 10:28:05 UTC 2022 with GitHub Actions
#
# Column D holds price text that looks numeric (e.g. "245.95", "0.05850").
# Assigning such a string straight to .Value lets Excel auto-coerce it to a
# real number (losing significant trailing zeros / formatting), so those
# cells are written with a leading apostrophe to force text, then restyled
# back to "Normal" so no stray quote-prefix / number-format style lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# --- Simple price (column D) updates that don't change row identity ---
Set-TextValue "D2"  "245.95"
Set-TextValue "D4"  "5.418"
Set-TextValue "D5"  "0.05850"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "6.333"
Set-TextValue "D8"  "0.8074"

# Row 9 (FTXToken) - price + volume(1h) label change
Set-TextValue "D9" "0.9691"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"

# --- Rows 10-18: the coin list rotated up by one (WazirX moved into row10's
# slot, ..., One wrapped around into row18), together with new price data ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1427"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07463"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03210"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03041"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.147"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09401"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001602"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04801"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005887"
$ws.Range("E18").Value = "17OneONE"

# --- Remaining simple price / label updates ---
Set-TextValue "D19" "0.006141"

Set-TextValue "D20" "0.004104"
$ws.Range("E20").Value = "19HotbitTokenHTB"

Set-TextValue "D21" "0.0009925"

Set-TextValue "D24" "2.228"

Set-TextValue "D26" "0.1295"

Set-TextValue "D27" "0.0003398"
$ws.Range("E27").Value = "26UpBotsUBXT"

Set-TextValue "D40" "0.03872"

Set-TextValue "D41" "0.006649"

Set-TextValue "D42" "0.1075"

Set-TextValue "D43" "0.002590"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.006389"

Set-TextValue "D45" "0.00005610"

Set-TextValue "D47" "0.3898"

Set-TextValue "D48" "0.1465"

Set-TextValue "D49" "0.00002099"

Set-TextValue "D50" "0.01009"
